# Fruta / hortaliza, semanal
# A new weekly price record (D=44476, i.e. 2021-10-07) is inserted for
# "Hortaliza, Terminal La Palmera de La Serena - Cilantro" right before the
# existing row that used to be row 75, pushing every row from the old 75
# down to 86 one row further (to 76..87). The new row re-uses the same
# constant columns (mercado, categoría, variedad, calidad, unidad, origen,
# clasificación, etc.) shared by every other record in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75 - shifts old rows 75:86 down to 76:87.
$ws.Rows("75:75").Insert()

$ws.Range("A75").Value = 8
$ws.Range("B75").Value = "Terminal La Palmera de La Serena"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44476
$ws.Range("E75").Value = 4
$ws.Range("F75").Value = 100112040
$ws.Range("G75").Value = "Cilantro"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 3080
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 2000
$ws.Range("M75").Value = 1750
$ws.Range("N75").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O75").Value = "Provincia del Elquí"
$ws.Range("P75").Value = 1167
$ws.Range("Q75").Value = 1.5
$ws.Range("R75").Value = "Hortaliza"
